# Corrections following third round of review
# Remove the "subgenus" column (header + its templated value) from the
# Materials worksheet. This shifts every following column one position
# to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$ws.Columns("AS:AS").Delete()

$wb.Save()
